# Apis development commited by Lakshmi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Policies table: fill in the previously-blank row 9 with a new "description" field
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "description"
$ws.Range("C9").Value = "Varchar"

# --- Customer table: rename existing row18 field from "age" to "email"
$ws.Range("B18").Value = "email"

# --- Insert a brand new row above row 19 (pushes old row19 "policyId" row down to row20)
$ws.Rows.Item(19).Insert()

# --- New row19: mobileNo field
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "mobileNo"
$ws.Range("C19").Value = "varchar"
$ws.Range("A19:D19").Style = $ws.Range("A18:D18").Style

# --- Row20 (previously row19, shifted down): bump the sequence number
$ws.Range("A20").Value = 8

# --- Update the sheet view to match the saved selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("D9").Select()
